# Auto update Excel log
# Appends new PIR/mmWave "Living Room" sensor events (rows 121-127) to the
# "mmWave" worksheet, extending the logged data from 2026-02-01 18:55 through 18:56.

$xlPasteValues = -4163

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New rows to append: Date, Timestamp, Hour, Location, Value, Status
$newRows = @(
    @("2026-02-01", "18:55:21", "18:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "18:55:26", "18:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "18:55:39", "18:00", "Living Room", "NO_MOTION_DETECTED", "Inactive"),
    @("2026-02-01", "18:55:49", "18:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "18:56:00", "18:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "18:56:10", "18:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "18:56:21", "18:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 121

# Helper cell used to stage text so that Excel's smart-input parser does not
# reinterpret date-like strings (e.g. "2026-02-01") as real dates. Assigning
# via a Formula that evaluates to text, then pasting the resulting value,
# keeps the destination cell's text intact and avoids leaving any stray
# number-format/style behind.
$helper = $ws.Cells.Item(500, 50)

$r = $startRow
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $cell = $ws.Cells.Item($r, $c)
        $helper.Formula = '="' + $val + '"'
        $helper.Copy()
        $cell.PasteSpecial($xlPasteValues)
        $c++
    }
    $r++
}

$helper.Clear()
$excel.CutCopyMode = $false
